$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 509.7618600982453
$ws.Range("B6").Value = 3397.25108424776
$ws.Range("B7").Value = 0.8909465299377146
$ws.Range("B8").Value = 0.8499487239957362
$ws.Range("B9").Value = 0.8218141097449367
$ws.Range("B10").Value = 0.8499487239957361
$ws.Range("B11").Value = 0.8704990043010362
$ws.Range("B12").Value = 0.08652364522839813
$ws.Range("B13").Value = 96.52039685939764
$ws.Range("B14").Value = 280.529027310251
$ws.Range("B15").Value = 1.739801570301179
